# Adjusting NOV13/14 observation locations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 / Row 12 observation-location swap:
#  - old row 11 (USACE / 01440)  -> becomes row 12
#  - old row 12 (USGS  / 291929089562600) -> becomes row 11, with the USGS
#    gage id corrected from 291929089562600 to 07380260
$ws.Range("A11").Value = "USGS"
$ws.Range("B11").Value = "07380260"
$ws.Range("A12").Value = "USACE"
$ws.Range("B12").Value = "01440"

# Column B was widened (best-fit) after the longer station ids were entered.
$ws.Columns.Item(2).ColumnWidth = 14.8

# The active cell/selection moved to E11 while reviewing the edit.
$ws.Range("E11").Select()
